$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.822.75"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.783.99"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.41"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.00"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").Value = "3.781.98"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.446"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  +4.72%  "
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.12"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "4.415.00"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.765.25"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "67.783.94"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.13"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "457.26"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.692"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.94"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  -5.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "3.928.24"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("E33").Value = "  -7.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.88"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.93"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0991"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.79"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.979"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.16"
$ws.Range("E41").Value = "  -7.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.08"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.10"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.02"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.36"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.27"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "384.40"
$ws.Range("E51").Value = "  -1.87%  "
